$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 7 (GAGE-1) flow/stage coordinates
$ws.Range("B7").Value = 427070.614
$ws.Range("C7").Value = 6654948.9349999996

# Add new row 8 for GAGE-2, reusing the coordinates previously held by row 7
$ws.Range("A8").Value = "GAGE-2"
$ws.Range("B8").Value = 427030.76
$ws.Range("C8").Value = 6654934.2599999998
$ws.Range("D8").Value = "INFLOW"

# Update selection to match the new active cell
$ws.Range("A8").Select() | Out-Null
